# Auto-generated edit script: apply profit-recalculation updates
# across multiple worksheets (per commit "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 175.09091
$ws.Cells.Item(33, 9).Value = 185.1
$ws.Cells.Item(33, 11).Value = 185.1
$ws.Cells.Item(33, 13).Value = 43.90000000000001

$ws.Cells.Item(88, 8).Value = 6917.7
$ws.Cells.Item(88, 9).Value = 1433
$ws.Cells.Item(88, 10).Value = 9268.286
$ws.Cells.Item(88, 11).Value = 1433
$ws.Cells.Item(88, 12).Value = 9268.286
$ws.Cells.Item(88, 13).Value = -1027
$ws.Cells.Item(88, 14).Value = -10080.286

$ws.Cells.Item(91, 8).Value = 6917.7
$ws.Cells.Item(91, 9).Value = 1433
$ws.Cells.Item(91, 10).Value = 9268.286
$ws.Cells.Item(91, 11).Value = 1433
$ws.Cells.Item(91, 12).Value = 9268.286
$ws.Cells.Item(91, 13).Value = -29
$ws.Cells.Item(91, 14).Value = -12076.286

$ws.Cells.Item(111, 8).Value = 1879
$ws.Cells.Item(111, 9).Value = 1879
$ws.Cells.Item(111, 11).Value = 5637
$ws.Cells.Item(111, 13).Value = -2570

$ws.Cells.Item(112, 8).Value = 1477.08
$ws.Cells.Item(112, 10).Value = 1477.08
$ws.Cells.Item(112, 12).Value = 4431.24
$ws.Cells.Item(112, 14).Value = -6647.24

$ws.Cells.Item(133, 8).Value = 59400
$ws.Cells.Item(133, 10).Value = 59400
$ws.Cells.Item(133, 12).Value = 59400
$ws.Cells.Item(133, 14).Value = -69520

$ws.Cells.Item(136, 8).Value = 70247.8
$ws.Cells.Item(136, 10).Value = 70247.8
$ws.Cells.Item(136, 12).Value = 70247.8
$ws.Cells.Item(136, 14).Value = -80447.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 293654.84
$ws.Cells.Item(2, 9).Value = 428520.06
$ws.Cells.Item(2, 11).Value = 428520.06
$ws.Cells.Item(2, 13).Value = -428407.06

$ws.Cells.Item(14, 8).Value = 16669167
$ws.Cells.Item(14, 10).Value = 3750
$ws.Cells.Item(14, 12).Value = 3750
$ws.Cells.Item(14, 14).Value = -4100

$ws.Cells.Item(32, 8).Value = 6791.4287
$ws.Cells.Item(32, 9).Value = 4370.6284
$ws.Cells.Item(32, 11).Value = 4370.6284
$ws.Cells.Item(32, 13).Value = -4083.6284

$ws.Cells.Item(45, 8).Value = 6430050
$ws.Cells.Item(45, 9).Value = 18000802
$ws.Cells.Item(45, 11).Value = 18000802
$ws.Cells.Item(45, 13).Value = -18000425

$ws.Cells.Item(61, 8).Value = 6024.48
$ws.Cells.Item(61, 9).Value = 6775.25
$ws.Cells.Item(61, 10).Value = 4689.778
$ws.Cells.Item(61, 11).Value = 6775.25
$ws.Cells.Item(61, 12).Value = 4689.778
$ws.Cells.Item(61, 13).Value = -6563.25
$ws.Cells.Item(61, 14).Value = -5113.778

$ws.Cells.Item(74, 8).Value = 1299.7307
$ws.Cells.Item(74, 9).Value = 483.27777
$ws.Cells.Item(74, 10).Value = 3136.75
$ws.Cells.Item(74, 11).Value = 483.27777
$ws.Cells.Item(74, 12).Value = 3136.75
$ws.Cells.Item(74, 13).Value = 390.72223
$ws.Cells.Item(74, 14).Value = -4884.75

$ws.Cells.Item(77, 8).Value = 1299.7307
$ws.Cells.Item(77, 9).Value = 483.27777
$ws.Cells.Item(77, 10).Value = 3136.75
$ws.Cells.Item(77, 11).Value = 2416.38885
$ws.Cells.Item(77, 12).Value = 15683.75
$ws.Cells.Item(77, 13).Value = 1951.61115
$ws.Cells.Item(77, 14).Value = -24419.75

$ws.Cells.Item(110, 8).Value = 927.2
$ws.Cells.Item(110, 9).Value = 250.66667
$ws.Cells.Item(110, 10).Value = 3633.3333
$ws.Cells.Item(110, 11).Value = 250.66667
$ws.Cells.Item(110, 12).Value = 3633.3333
$ws.Cells.Item(110, 13).Value = 1794.33333
$ws.Cells.Item(110, 14).Value = -7723.3333

$ws.Cells.Item(116, 8).Value = 293654.84
$ws.Cells.Item(116, 9).Value = 428520.06
$ws.Cells.Item(116, 11).Value = 428520.06
$ws.Cells.Item(116, 13).Value = -426226.06

$ws.Cells.Item(122, 8).Value = 1359
$ws.Cells.Item(122, 9).Value = 931.75
$ws.Cells.Item(122, 10).Value = 2042.6
$ws.Cells.Item(122, 11).Value = 2795.25
$ws.Cells.Item(122, 12).Value = 6127.799999999999
$ws.Cells.Item(122, 13).Value = -345.25
$ws.Cells.Item(122, 14).Value = -11027.8

$ws.Cells.Item(132, 8).Value = 1923.125
$ws.Cells.Item(132, 9).Value = 1479.1177
$ws.Cells.Item(132, 11).Value = 4437.3531
$ws.Cells.Item(132, 13).Value = -1907.3531

$ws.Cells.Item(136, 8).Value = 6024.48
$ws.Cells.Item(136, 9).Value = 6775.25
$ws.Cells.Item(136, 10).Value = 4689.778
$ws.Cells.Item(136, 11).Value = 20325.75
$ws.Cells.Item(136, 12).Value = 14069.334
$ws.Cells.Item(136, 13).Value = -17775.75
$ws.Cells.Item(136, 14).Value = -19169.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 293654.84
$ws.Cells.Item(3, 9).Value = 428520.06
$ws.Cells.Item(3, 11).Value = 428520.06
$ws.Cells.Item(3, 13).Value = -428406.06

$ws.Cells.Item(22, 8).Value = 893.7778
$ws.Cells.Item(22, 9).Value = 837.8
$ws.Cells.Item(22, 11).Value = 837.8
$ws.Cells.Item(22, 13).Value = -664.8

$ws.Cells.Item(134, 8).Value = 5410.4644
$ws.Cells.Item(134, 9).Value = 6153.1816
$ws.Cells.Item(134, 10).Value = 2687.1667
$ws.Cells.Item(134, 11).Value = 18459.5448
$ws.Cells.Item(134, 12).Value = 8061.500100000001
$ws.Cells.Item(134, 13).Value = -15924.5448
$ws.Cells.Item(134, 14).Value = -13131.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1540
$ws.Cells.Item(22, 9).Value = 666.6667
$ws.Cells.Item(22, 10).Value = 1914.2858
$ws.Cells.Item(22, 11).Value = 666.6667
$ws.Cells.Item(22, 12).Value = 1914.2858
$ws.Cells.Item(22, 13).Value = -316.6667
$ws.Cells.Item(22, 14).Value = -2614.2858

$ws.Cells.Item(31, 8).Value = 2272.2222
$ws.Cells.Item(31, 9).Value = 2125
$ws.Cells.Item(31, 11).Value = 2125
$ws.Cells.Item(31, 13).Value = -1830

$ws.Cells.Item(34, 8).Value = 2272.2222
$ws.Cells.Item(34, 9).Value = 2125
$ws.Cells.Item(34, 11).Value = 2125
$ws.Cells.Item(34, 13).Value = -1923

$ws.Cells.Item(92, 8).Value = 39999
$ws.Cells.Item(92, 10).Value = 39999
$ws.Cells.Item(92, 12).Value = 39999
$ws.Cells.Item(92, 14).Value = -44991

$ws.Cells.Item(105, 8).Value = 829.2
$ws.Cells.Item(105, 9).Value = 797.8889
$ws.Cells.Item(105, 11).Value = 797.8889
$ws.Cells.Item(105, 13).Value = 949.1111

$ws.Cells.Item(107, 8).Value = 1088.4286
$ws.Cells.Item(107, 10).Value = 1171.1428
$ws.Cells.Item(107, 12).Value = 1171.1428
$ws.Cells.Item(107, 14).Value = -5011.1428

$ws.Cells.Item(122, 8).Value = 2975.6
$ws.Cells.Item(122, 10).Value = 6671
$ws.Cells.Item(122, 12).Value = 20013
$ws.Cells.Item(122, 14).Value = -24913

$ws.Cells.Item(132, 8).Value = 2268.28
$ws.Cells.Item(132, 9).Value = 1417.5883
$ws.Cells.Item(132, 11).Value = 4252.7649
$ws.Cells.Item(132, 13).Value = -1722.7649

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 169.8
$ws.Cells.Item(10, 9).Value = 169.8
$ws.Cells.Item(10, 11).Value = 509.4
$ws.Cells.Item(10, 13).Value = -370.4

$ws.Cells.Item(87, 8).Value = 188239
$ws.Cells.Item(87, 9).Value = 355228
$ws.Cells.Item(87, 10).Value = 21250
$ws.Cells.Item(87, 11).Value = 1065684
$ws.Cells.Item(87, 12).Value = 63750
$ws.Cells.Item(87, 13).Value = -1064436
$ws.Cells.Item(87, 14).Value = -66246

$ws.Cells.Item(90, 8).Value = 188239
$ws.Cells.Item(90, 9).Value = 355228
$ws.Cells.Item(90, 10).Value = 21250
$ws.Cells.Item(90, 11).Value = 3197052
$ws.Cells.Item(90, 12).Value = 191250
$ws.Cells.Item(90, 13).Value = -3190812
$ws.Cells.Item(90, 14).Value = -203730

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1400.6666
$ws.Cells.Item(113, 10).Value = 1522
$ws.Cells.Item(113, 12).Value = 1522
$ws.Cells.Item(113, 14).Value = -5862

$ws.Cells.Item(132, 8).Value = 2265178.2
$ws.Cells.Item(132, 9).Value = 3206844.8
$ws.Cells.Item(132, 11).Value = 9620534.399999999
$ws.Cells.Item(132, 13).Value = -9618004.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2334.0557
$ws.Cells.Item(61, 9).Value = 2179.4546
$ws.Cells.Item(61, 10).Value = 2577
$ws.Cells.Item(61, 11).Value = 2179.4546
$ws.Cells.Item(61, 12).Value = 2577
$ws.Cells.Item(61, 13).Value = -1977.4546
$ws.Cells.Item(61, 14).Value = -2981

$ws.Cells.Item(113, 8).Value = 2334.0557
$ws.Cells.Item(113, 9).Value = 2179.4546
$ws.Cells.Item(113, 10).Value = 2577
$ws.Cells.Item(113, 11).Value = 2179.4546
$ws.Cells.Item(113, 12).Value = 2577
$ws.Cells.Item(113, 13).Value = -9.454600000000028
$ws.Cells.Item(113, 14).Value = -6917

$ws.Cells.Item(132, 8).Value = 1539.1136
$ws.Cells.Item(132, 9).Value = 1202.0476
$ws.Cells.Item(132, 10).Value = 1846.8695
$ws.Cells.Item(132, 11).Value = 3606.142800000001
$ws.Cells.Item(132, 12).Value = 5540.6085
$ws.Cells.Item(132, 13).Value = -1076.142800000001
$ws.Cells.Item(132, 14).Value = -10600.6085

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1190.2188
$ws.Cells.Item(132, 9).Value = 878.0851
$ws.Cells.Item(132, 10).Value = 2053.1765
$ws.Cells.Item(132, 11).Value = 2634.2553
$ws.Cells.Item(132, 12).Value = 6159.529500000001
$ws.Cells.Item(132, 13).Value = -104.2552999999998
$ws.Cells.Item(132, 14).Value = -11219.5295
